$d = $word.ActiveDocument

# Locate the whole final run of the paragraph describing write_file_sort_str(...): the
# sentence "который нужно прочесть, ... метод - числовой или строковый. ". We match the
# complete run text so replacement keeps it as a single, unmodified run (matching the
# diff, which leaves this run untouched and only appends new runs after it).
$dash = [char]0x2013
$anchor = "который нужно прочесть, путь файла, в который нужно записать данные, параметр, по которому требуется сортировка и метод $dash числовой или строковый. "

$r = $d.Content
$found = $r.Find.Execute($anchor, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find anchor sentence"
}

# Work off of a brand new Range built from the match boundaries (re-using the Find-mutated
# range object directly for XML/formatting operations is unreliable in this host).
$target = $d.Range($r.Start, $r.End)

$xml = @'
<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r w:rsidR="006D5070"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve">который нужно прочесть, путь файла, в который нужно записать данные, параметр, по которому требуется сортировка и метод &#8211; числовой или строковый. </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve"> Также написана функция для считывания кол-ва файлов в директории </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>countFilesInDir</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>(</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr><w:t>path</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve">). </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$target.InsertXML($xml) | Out-Null
